# Auto-generated edit script: updates Final Fantasy XIV market-board price/profit
# figures across the 8 crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to match a refreshed Universalis price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 70000
$ws.Range("J3").Value = 70000
$ws.Range("L3").Value = 70000
$ws.Range("N3").Value = -70228

$ws.Range("H70").Value = 54409.43
$ws.Range("J70").Value = 89554.414
$ws.Range("L70").Value = 268663.242
$ws.Range("N70").Value = -269203.242

$ws.Range("H73").Value = 54409.43
$ws.Range("J73").Value = 89554.414
$ws.Range("L73").Value = 268663.242
$ws.Range("N73").Value = -270535.242

$ws.Range("H86").Value = 1786.55
$ws.Range("I86").Value = 1651.7142
$ws.Range("J86").Value = 2101.1667
$ws.Range("K86").Value = 1651.7142
$ws.Range("L86").Value = 2101.1667
$ws.Range("M86").Value = -528.7141999999999
$ws.Range("N86").Value = -4347.1667

$ws.Range("H89").Value = 1786.55
$ws.Range("I89").Value = 1651.7142
$ws.Range("J89").Value = 2101.1667
$ws.Range("K89").Value = 8258.571
$ws.Range("L89").Value = 10505.8335
$ws.Range("M89").Value = -2642.571
$ws.Range("N89").Value = -21737.8335

$ws.Range("H98").Value = 2600.0625
$ws.Range("I98").Value = 2600.0625
$ws.Range("K98").Value = 2600.0625
$ws.Range("M98").Value = -1102.0625

$ws.Range("H102").Value = 70000
$ws.Range("J102").Value = 70000
$ws.Range("L102").Value = 70000
$ws.Range("N102").Value = -76490

$ws.Range("H122").Value = 2600.0625
$ws.Range("I122").Value = 2600.0625
$ws.Range("K122").Value = 7800.1875
$ws.Range("M122").Value = -5350.1875

$ws.Range("H138").Value = 45556.74
$ws.Range("I138").Value = 2075.6667
$ws.Range("J138").Value = 502108
$ws.Range("K138").Value = 6227.000100000001
$ws.Range("L138").Value = 1506324
$ws.Range("M138").Value = -1087.000100000001
$ws.Range("N138").Value = -1516604

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4705.5
$ws.Range("I88").Value = 1318.75
$ws.Range("J88").Value = 6060.2
$ws.Range("K88").Value = 1318.75
$ws.Range("L88").Value = 6060.2
$ws.Range("M88").Value = -912.75
$ws.Range("N88").Value = -6872.2

$ws.Range("H91").Value = 4705.5
$ws.Range("I91").Value = 1318.75
$ws.Range("J91").Value = 6060.2
$ws.Range("K91").Value = 1318.75
$ws.Range("L91").Value = 6060.2
$ws.Range("M91").Value = 85.25
$ws.Range("N91").Value = -8868.200000000001

$ws.Range("H95").Value = 302219.8
$ws.Range("J95").Value = 302219.8
$ws.Range("L95").Value = 302219.8
$ws.Range("N95").Value = -307711.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 17828.5
$ws.Range("J103").Value = 17828.5
$ws.Range("L103").Value = 17828.5
$ws.Range("N103").Value = -20172.5

$ws.Range("H134").Value = 2526.3215
$ws.Range("J134").Value = 4507.875
$ws.Range("L134").Value = 13523.625
$ws.Range("N134").Value = -18593.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13275.593
$ws.Range("I58").Value = 1405.72
$ws.Range("J58").Value = 161649
$ws.Range("K58").Value = 1405.72
$ws.Range("L58").Value = 161649
$ws.Range("M58").Value = -1202.72
$ws.Range("N58").Value = -162055

$ws.Range("H132").Value = 39552.27
$ws.Range("I132").Value = 46366.316
$ws.Range("K132").Value = 139098.948
$ws.Range("M132").Value = -136568.948

$ws.Range("H136").Value = 13275.593
$ws.Range("I136").Value = 1405.72
$ws.Range("J136").Value = 161649
$ws.Range("K136").Value = 4217.16
$ws.Range("L136").Value = 484947
$ws.Range("M136").Value = -1667.16
$ws.Range("N136").Value = -490047

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 309.5
$ws.Range("I18").Value = 331.6
$ws.Range("K18").Value = 994.8000000000001
$ws.Range("M18").Value = -825.8000000000001

$ws.Range("H39").Value = 8637.5
$ws.Range("J39").Value = 8637.5
$ws.Range("L39").Value = 25912.5
$ws.Range("N39").Value = -26500.5

$ws.Range("H45").Value = 2278.1
$ws.Range("I45").Value = 1533
$ws.Range("J45").Value = 2597.4285
$ws.Range("K45").Value = 4599
$ws.Range("L45").Value = 7792.2855
$ws.Range("M45").Value = -4067
$ws.Range("N45").Value = -8856.2855

$ws.Range("H51").Value = 4391.304
$ws.Range("I51").Value = 2200
$ws.Range("K51").Value = 6600
$ws.Range("M51").Value = -6140

$ws.Range("H81").Value = 4707.362
$ws.Range("I81").Value = 2811.5
$ws.Range("K81").Value = 8434.5
$ws.Range("M81").Value = -7311.5

$ws.Range("H84").Value = 4707.362
$ws.Range("I84").Value = 2811.5
$ws.Range("K84").Value = 25303.5
$ws.Range("M84").Value = -19687.5

$ws.Range("H102").Value = 4341
$ws.Range("J102").Value = 3999
$ws.Range("L102").Value = 11997
$ws.Range("N102").Value = -16865

$ws.Range("H103").Value = 751.94446
$ws.Range("I103").Value = 421.75
$ws.Range("J103").Value = 1016.1
$ws.Range("K103").Value = 1265.25
$ws.Range("L103").Value = 3048.3
$ws.Range("M103").Value = -386.25
$ws.Range("N103").Value = -4806.3

$ws.Range("H104").Value = 1575
$ws.Range("J104").Value = 950
$ws.Range("L104").Value = 2850
$ws.Range("N104").Value = -8092

$ws.Range("H112").Value = 3212.5
$ws.Range("I112").Value = 2350
$ws.Range("K112").Value = 7050
$ws.Range("M112").Value = -5942

$ws.Range("H116").Value = 9459.405000000001
$ws.Range("I116").Value = 3332.6667
$ws.Range("K116").Value = 9998.000100000001
$ws.Range("M116").Value = -6556.000100000001

$ws.Range("H118").Value = 1275.8
$ws.Range("I118").Value = 844.75
$ws.Range("K118").Value = 2534.25
$ws.Range("M118").Value = -1291.25

$ws.Range("H119").Value = 1069.75
$ws.Range("I119").Value = 864.5
$ws.Range("J119").Value = 1275
$ws.Range("K119").Value = 2593.5
$ws.Range("L119").Value = 3825
$ws.Range("M119").Value = 2244.5
$ws.Range("N119").Value = -13501

$ws.Range("H120").Value = 20000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 60000
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -69676

$ws.Range("H121").Value = 67996.81
$ws.Range("I121").Value = 213584
$ws.Range("J121").Value = 1820.8182
$ws.Range("K121").Value = 640752
$ws.Range("L121").Value = 5462.4546
$ws.Range("M121").Value = -639442
$ws.Range("N121").Value = -8082.4546

$ws.Range("H124").Value = 4312
$ws.Range("J124").Value = 4871.4287
$ws.Range("L124").Value = 14614.2861
$ws.Range("N124").Value = -24434.2861

$ws.Range("H137").Value = 3194.7058
$ws.Range("J137").Value = 4063.2
$ws.Range("L137").Value = 12189.6
$ws.Range("N137").Value = -22389.6

$ws.Range("H139").Value = 1617.4
$ws.Range("I139").Value = 1259.8
$ws.Range("J139").Value = 1975
$ws.Range("K139").Value = 3779.4
$ws.Range("L139").Value = 5925
$ws.Range("M139").Value = 1360.6
$ws.Range("N139").Value = -16205

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 2458.5
$ws.Range("J92").Value = 2458.5
$ws.Range("L92").Value = 2458.5
$ws.Range("N92").Value = -6202.5

$ws.Range("H98").Value = 29500
$ws.Range("J98").Value = 29500
$ws.Range("L98").Value = 29500
$ws.Range("N98").Value = -35490

$ws.Range("H105").Value = 2539917.8
$ws.Range("J105").Value = 3366557
$ws.Range("L105").Value = 3366557
$ws.Range("N105").Value = -3373545

$ws.Range("H126").Value = 2766
$ws.Range("I126").Value = 2127.7144
$ws.Range("K126").Value = 6383.1432
$ws.Range("M126").Value = -3913.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 8024.75
$ws.Range("J106").Value = 8024.75
$ws.Range("L106").Value = 8024.75
$ws.Range("N106").Value = -10548.75

$ws.Range("H132").Value = 2102.8
$ws.Range("I132").Value = 941.75
$ws.Range("J132").Value = 3429.7144
$ws.Range("K132").Value = 2825.25
$ws.Range("L132").Value = 10289.1432
$ws.Range("M132").Value = -295.25
$ws.Range("N132").Value = -15349.1432

$ws.Range("H136").Value = 2938.8206
$ws.Range("I136").Value = 2845.4827
$ws.Range("J136").Value = 3209.5
$ws.Range("K136").Value = 8536.4481
$ws.Range("L136").Value = 9628.5
$ws.Range("M136").Value = -5986.4481
$ws.Range("N136").Value = -14728.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 179515.08
$ws.Range("I122").Value = 288962.12
$ws.Range("K122").Value = 866886.36
$ws.Range("M122").Value = -864436.36

$ws.Range("H132").Value = 1621.0416
$ws.Range("I132").Value = 1350.1818
$ws.Range("K132").Value = 4050.5454
$ws.Range("M132").Value = -1520.5454
